# "Fruta / hortaliza, semanal" weekly data update:
# insert one new weekly record (a new row) for Choclo / Dulce o Americano / Primera
# at row 1134, pushing the subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 1134; existing rows 1134:1193 shift down to 1135:1194.
$ws.Rows.Item(1134).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(1134, 1).Value  = 6
$ws.Cells.Item(1134, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1134, 3).Value  = "Metropolitana"
$ws.Cells.Item(1134, 4).Value  = 44753
$ws.Cells.Item(1134, 5).Value  = 13
$ws.Cells.Item(1134, 6).Value  = 100112024
$ws.Cells.Item(1134, 7).Value  = "Choclo"
$ws.Cells.Item(1134, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(1134, 9).Value  = "Primera"
$ws.Cells.Item(1134, 10).Value = 400
$ws.Cells.Item(1134, 11).Value = 35000
$ws.Cells.Item(1134, 12).Value = 37000
$ws.Cells.Item(1134, 13).Value = 35850
$ws.Cells.Item(1134, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(1134, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1134, 16).Value = 512
$ws.Cells.Item(1134, 17).Value = 70
$ws.Cells.Item(1134, 18).Value = "Hortaliza"
